# Update of all values to match PDF edition 10 (commit 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update the "Total" label to "Total/average" (shared string used by A8/A15)
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Total/average"
$ws.Range("A15").Value = "Total/average"

# ---------------------------------------------------------------------------
# 2. Update the data values (rows 2-15, columns B:E)
# ---------------------------------------------------------------------------
$data = @{
    2  = @(2028273, 466, 97, 3144)
    3  = @(1423613, 706, 127, 4695)
    4  = @(595555, 402, 92, 861)
    5  = @(209977, 586, 118, 3538)
    6  = @(265609, 189, 88, 392)
    7  = @(229566, 478, 99, 5955)
    8  = @(4752593, 471, 104, 3098)
    9  = @(671271, 2855, 396, 36909)
    10 = @(97720, 1344, 212, 9106)
    11 = @(87016, 1969, 291, 5856)
    12 = @(151769, 1933, 289, 19996)
    13 = @(8142, 2331, 343, 23044)
    14 = @(113382, 3310, 450, 60065)
    15 = @(1129300, 2290, 330, 25829)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
}

# ---------------------------------------------------------------------------
# 3. Bold the "Total/average" summary rows (8 and 15) across A:E
# ---------------------------------------------------------------------------
$ws.Range("A8:E8").Font.Bold = $true
$ws.Range("A15:E15").Font.Bold = $true

# ---------------------------------------------------------------------------
# 4. Update the active selection on the sheet view
# ---------------------------------------------------------------------------
$ws.Range("G12").Select()
